$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header text (Volume/Number and report week dates) ---
$ws.Range("A8").Value = "Volume 33   Number  3"
$ws.Range("C9").Value = "Report Covering the Week  1/12/2026  Through  1/18/2026"

# --- Phase 1: fix up cell styles (number format) / text-placeholder cells ---
# Uses Copy($dest) (format+value) from stable template cells so the destination
# lands on the exact same shared cell-style index as the target workbook.
# Row 23 is untouched by this edit, so it is used as the stable source for the
# "no data" placeholder cells (shared strings "0" / "***.*").
$ws.Range("C15").Copy($ws.Range("D14"))
$ws.Range("M15").Copy($ws.Range("E14"))
$ws.Range("C15").Copy($ws.Range("G14"))
$ws.Range("M15").Copy($ws.Range("H14"))
$ws.Range("C15").Copy($ws.Range("J14"))
$ws.Range("M15").Copy($ws.Range("K14"))
$ws.Range("M15").Copy($ws.Range("N14"))
$ws.Range("C15").Copy($ws.Range("D15"))
$ws.Range("M15").Copy($ws.Range("E15"))
$ws.Range("C15").Copy($ws.Range("G15"))
$ws.Range("M15").Copy($ws.Range("H15"))
$ws.Range("C15").Copy($ws.Range("J15"))
$ws.Range("M15").Copy($ws.Range("K15"))
$ws.Range("C23").Copy($ws.Range("C16"))
$ws.Range("C15").Copy($ws.Range("D22"))
$ws.Range("M15").Copy($ws.Range("E22"))
$ws.Range("C15").Copy($ws.Range("J22"))
$ws.Range("M15").Copy($ws.Range("K22"))
$ws.Range("C15").Copy($ws.Range("D27"))
$ws.Range("M15").Copy($ws.Range("E27"))
$ws.Range("C23").Copy($ws.Range("D28"))
$ws.Range("E23").Copy($ws.Range("E28"))

# --- Phase 2: write final numeric values ---
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = -100
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = -100
$ws.Range("J14").Value = 1
$ws.Range("K14").Value = -100
$ws.Range("N14").Value = -100
$ws.Range("C15").Value = 2
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 100
$ws.Range("F15").Value = 3
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 200
$ws.Range("I15").Value = 3
$ws.Range("J15").Value = 1
$ws.Range("K15").Value = 200
$ws.Range("M15").Value = 200
$ws.Range("N15").Value = 200
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -100
$ws.Range("F16").Value = 13
$ws.Range("H16").Value = 44.444444444444
$ws.Range("J16").Value = 7
$ws.Range("K16").Value = 14.285714285714
$ws.Range("L16").Value = -42.857142857142
$ws.Range("M16").Value = -50
$ws.Range("N16").Value = -89.473684210526
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = -50
$ws.Range("F17").Value = 10
$ws.Range("G17").Value = 11
$ws.Range("H17").Value = -9.090909090909
$ws.Range("I17").Value = 7
$ws.Range("J17").Value = 8
$ws.Range("K17").Value = -12.5
$ws.Range("L17").Value = -36.363636363636
$ws.Range("M17").Value = 16.666666666666
$ws.Range("N17").Value = -61.111111111111
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -75
$ws.Range("F18").Value = 22
$ws.Range("G18").Value = 25
$ws.Range("H18").Value = -12
$ws.Range("I18").Value = 10
$ws.Range("J18").Value = 12
$ws.Range("K18").Value = -16.666666666666
$ws.Range("L18").Value = -23.076923076923
$ws.Range("M18").Value = 42.857142857142
$ws.Range("N18").Value = -86.301369863013
$ws.Range("C19").Value = 19
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = 46.153846153846
$ws.Range("F19").Value = 54
$ws.Range("G19").Value = 42
$ws.Range("H19").Value = 28.571428571428
$ws.Range("I19").Value = 39
$ws.Range("J19").Value = 22
$ws.Range("K19").Value = 77.272727272727
$ws.Range("L19").Value = -18.75
$ws.Range("M19").Value = 105.263157894737
$ws.Range("N19").Value = 5.405405405405
$ws.Range("C20").Value = 10
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 400
$ws.Range("F20").Value = 27
$ws.Range("G20").Value = 11
$ws.Range("H20").Value = 145.454545454545
$ws.Range("I20").Value = 21
$ws.Range("J20").Value = 7
$ws.Range("K20").Value = 200
$ws.Range("L20").Value = 425
$ws.Range("M20").Value = 162.5
$ws.Range("N20").Value = -82.5
$ws.Range("C21").Value = 33
$ws.Range("D21").Value = 27
$ws.Range("E21").Value = 22.222222222222
$ws.Range("F21").Value = 129
$ws.Range("G21").Value = 100
$ws.Range("H21").Value = 29
$ws.Range("I21").Value = 88
$ws.Range("J21").Value = 58
$ws.Range("K21").Value = 51.724137931034
$ws.Range("L21").Value = -2.222222222222
$ws.Range("M21").Value = 54.38596491228
$ws.Range("N21").Value = -73.006134969325
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 0
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = 66.666666666666
$ws.Range("J22").Value = 1
$ws.Range("K22").Value = 200
$ws.Range("M22").Value = 0
$ws.Range("C24").Value = 35
$ws.Range("D24").Value = 29
$ws.Range("E24").Value = 20.689655172413
$ws.Range("F24").Value = 114
$ws.Range("G24").Value = 130
$ws.Range("H24").Value = -12.307692307692
$ws.Range("I24").Value = 72
$ws.Range("J24").Value = 90
$ws.Range("K24").Value = -20
$ws.Range("L24").Value = -21.739130434782
$ws.Range("M24").Value = 125
$ws.Range("C25").Value = 22
$ws.Range("D25").Value = 18
$ws.Range("E25").Value = 22.222222222222
$ws.Range("F25").Value = 67
$ws.Range("G25").Value = 73
$ws.Range("H25").Value = -8.219178082191
$ws.Range("I25").Value = 42
$ws.Range("J25").Value = 56
$ws.Range("K25").Value = -25
$ws.Range("L25").Value = -23.636363636363
$ws.Range("C26").Value = 12
$ws.Range("D26").Value = 11
$ws.Range("E26").Value = 9.090909090909
$ws.Range("G26").Value = 33
$ws.Range("H26").Value = 6.060606060606
$ws.Range("I26").Value = 24
$ws.Range("J26").Value = 27
$ws.Range("K26").Value = -11.111111111111
$ws.Range("L26").Value = 9.090909090909
$ws.Range("M26").Value = 41.176470588235
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 100
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 50
$ws.Range("I27").Value = 3
$ws.Range("J27").Value = 2
$ws.Range("K27").Value = 50
$ws.Range("L27").Value = 200
$ws.Range("F28").Value = 3
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 50
